# Add overview tables for omics analysis tools
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data for the omics tools tables (row 3 cells already existed
# with formatting; rows 4-5 gain brand new A/B cells).
$ws.Range("A3").Value = "t2_omics_analysis_tools"
$ws.Range("B3").Value = "Analysis tools for omics data"

$ws.Range("A4").Value = "t3_omics_post-analysis_tools"
$ws.Range("B4").Value = "Post-analysis tools for omics data"

$ws.Range("A5").Value = "t4_gene_enrichment_tools"
$ws.Range("B5").Value = "Gene enrichment analysis tools"

# Widen column B to fit the longer description text (target stored width
# 40.98 chars; the engine snaps ColumnWidth to a pixel grid on save, so
# 40.15 is the input that lands on the closest achievable value, 41).
$ws.Columns.Item(2).ColumnWidth = 40.15

# Move the active selection to A3.
$ws.Range("A3").Select() | Out-Null
